# "fixed export and fixing maps"
# The sheet previously compared Area across three census years
# (1989 / 2002 / 2014). This reverts it back to a simple, single-year
# (2014) export: drop the "(according to the population census data)"
# caption, drop the 1989 and 2002 columns, and pad the table out with a
# few blank, taller rows underneath (matching the simpler export layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "(according to the population census data)" caption in A2 -
# the row stays, just empty now.
$ws.Range("A2").Clear()

# The blank spacer row (old row 3) disappears entirely; everything below
# shifts up by one.
$ws.Rows("3:3").Delete()

# Only the 2014 figures are kept; 1989 (old col B) and 2002 (old col C)
# are removed, so the old "2014" column becomes column B.
$ws.Columns("B:C").Delete()

# The remaining rows get a taller, consistent row height, and a handful
# of extra blank (but equally tall) rows are left under the table.
$ws.Rows("1:11").RowHeight = 20.1
